{"js": "// The worksheet/table holds 5 \"content\" rows (0, 4, 8, 12, 16 \u2014 the rows in\n// between are intentionally blank spacer rows) x 5 columns of\n// \"A\u00f7B=C, D\" style division answers. Replace each cell's text in place,\n// keyed by its (row, column) position, so we don't have to worry about any\n// value collisions between the old and new text (some new answers equal an\n// old answer found elsewhere in the table).\nconst table = context.document.body.tables.getFirst();\n\nconst replacements = [\n  // row 0\n  { row: 0, col: 0, text: \"63\u00f77=9, 0\" },\n  { row: 0, col: 1, text: \"79\u00f77=11, 2\" },\n  { row: 0, col: 2, text: \"39\u00f75=7, 4\" },\n  { row: 0, col: 3, text: \"51\u00f74=12, 3\" },\n  { row: 0, col: 4, text: \"47\u00f76=7, 5\" },\n  // row 4\n  { row: 4, col: 0, text: \"93\u00f78=11, 5\" },\n  { row: 4, col: 1, text: \"48\u00f77=6, 6\" },\n  { row: 4, col: 2, text: \"27\u00f78=3, 3\" },\n  { row: 4, col: 3, text: \"18\u00f73=6, 0\" },\n  { row: 4, col: 4, text: \"47\u00f73=15, 2\" },\n  // row 8\n  { row: 8, col: 0, text: \"16\u00f76=2, 4\" },\n  { row: 8, col: 1, text: \"65\u00f78=8, 1\" },\n  { row: 8, col: 2, text: \"31\u00f72=15, 1\" },\n  { row: 8, col: 3, text: \"69\u00f73=23, 0\" },\n  { row: 8, col: 4, text: \"49\u00f74=12, 1\" },\n  // row 12\n  { row: 12, col: 0, text: \"80\u00f73=26, 2\" },\n  { row: 12, col: 1, text: \"14\u00f72=7, 0\" },\n  { row: 12, col: 2, text: \"99\u00f74=24, 3\" },\n  { row: 12, col: 3, text: \"26\u00f72=13, 0\" },\n  { row: 12, col: 4, text: \"53\u00f73=17, 2\" },\n  // row 16\n  { row: 16, col: 0, text: \"26\u00f75=5, 1\" },\n  { row: 16, col: 1, text: \"11\u00f78=1, 3\" },\n  { row: 16, col: 2, text: \"68\u00f74=17, 0\" },\n  { row: 16, col: 3, text: \"64\u00f78=8, 0\" },\n  { row: 16, col: 4, text: \"21\u00f78=2, 5\" },\n];\n\nfor (const { row, col, text } of replacements) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# The table has 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17 in\n# Word's 1-based COM indexing) actually holds an \"A\u00f7B=C, D\" division-answer\n# string -- the rows in between are blank spacer rows. Update each populated\n# cell's text in place by (row, column) position so we don't have to worry\n# about collisions between old/new values (some new answers duplicate an old\n# answer found elsewhere in the table).\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$t.Cell(1, 1).Range.Text  = \"63\u00f77=9, 0\"\n$t.Cell(1, 2).Range.Text  = \"79\u00f77=11, 2\"\n$t.Cell(1, 3).Range.Text  = \"39\u00f75=7, 4\"\n$t.Cell(1, 4).Range.Text  = \"51\u00f74=12, 3\"\n$t.Cell(1, 5).Range.Text  = \"47\u00f76=7, 5\"\n\n$t.Cell(5, 1).Range.Text  = \"93\u00f78=11, 5\"\n$t.Cell(5, 2).Range.Text  = \"48\u00f77=6, 6\"\n$t.Cell(5, 3).Range.Text  = \"27\u00f78=3, 3\"\n$t.Cell(5, 4).Range.Text  = \"18\u00f73=6, 0\"\n$t.Cell(5, 5).Range.Text  = \"47\u00f73=15, 2\"\n\n$t.Cell(9, 1).Range.Text  = \"16\u00f76=2, 4\"\n$t.Cell(9, 2).Range.Text  = \"65\u00f78=8, 1\"\n$t.Cell(9, 3).Range.Text  = \"31\u00f72=15, 1\"\n$t.Cell(9, 4).Range.Text  = \"69\u00f73=23, 0\"\n$t.Cell(9, 5).Range.Text  = \"49\u00f74=12, 1\"\n\n$t.Cell(13, 1).Range.Text = \"80\u00f73=26, 2\"\n$t.Cell(13, 2).Range.Text = \"14\u00f72=7, 0\"\n$t.Cell(13, 3).Range.Text = \"99\u00f74=24, 3\"\n$t.Cell(13, 4).Range.Text = \"26\u00f72=13, 0\"\n$t.Cell(13, 5).Range.Text = \"53\u00f73=17, 2\"\n\n$t.Cell(17, 1).Range.Text = \"26\u00f75=5, 1\"\n$t.Cell(17, 2).Range.Text = \"11\u00f78=1, 3\"\n$t.Cell(17, 3).Range.Text = \"68\u00f74=17, 0\"\n$t.Cell(17, 4).Range.Text = \"64\u00f78=8, 0\"\n$t.Cell(17, 5).Range.Text = \"21\u00f78=2, 5\"\n"}
